# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values computed for each game row (2-16),
# replacing the previous "Strike#" values in column G.
$kValues = @{
    2  = 3
    3  = 3
    4  = 4
    5  = 2
    6  = 2
    7  = 3
    8  = 4
    9  = 7
    10 = 6
    11 = 4
    12 = 1
    13 = 1
    14 = 4
    15 = 1
    16 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
